# Weekly price update: insert a new record row for
# "Macroferia Regional de Talca" - Repollo (Hortaliza) right before the
# existing row 135, shifting all following rows (and their data) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(135).Insert()

$ws.Range("A135").Value = 5
$ws.Range("B135").Value = "Macroferia Regional de Talca"
$ws.Range("C135").Value = "Maule"
$ws.Range("D135").Value2 = 44468
$ws.Range("E135").Value = 7
$ws.Range("F135").Value = 100112006
$ws.Range("G135").Value = "Repollo"
$ws.Range("H135").Value = "Crespo record"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 450
$ws.Range("L135").Value = 450
$ws.Range("M135").Value = 450
$ws.Range("N135").Value = "$/unidad"
$ws.Range("O135").Value = "Región del Maule"
$ws.Range("P135").Value = 450
$ws.Range("Q135").Value = 1
$ws.Range("R135").Value = "Hortaliza"

Write-Host "Inserted new weekly Repollo record at row 135"
